# Removed Test Case Inter-Dependency
#
# - ProductLoanInput!B1 (productname) and ProductLoanOutput!B1 (verifyloanproduct
#   readback) get an "-1st" suffix appended to the product name so this test
#   case no longer collides with another run.
# - ProductLoanInput!B2 (shortname) switches from the numeric default 2470 to
#   the distinct text value "247d".
# - Selection/active-tab bookkeeping moves from ProductLoanInput (cell B8) to
#   ProductLoanOutput being the active sheet, each sheet's selection resets to
#   B1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2470-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-CASH-1st"

# Product name (shared between the input sheet and the output readback sheet)
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Short name becomes a distinct text value instead of the numeric default
$ws1.Range("B2").Value = "247d"

# Nudge the output cell's font so it re-resolves to the plain (non-bold-family)
# cell style used elsewhere in the sheet.
$ws2.Range("B1").Font.Name = "Arial"

# Selection / active sheet bookkeeping
$null = $ws1.Range("B1").Select()
$null = $ws2.Activate()
